# Fixed int/float bug when reading data from ByBit in find_crossing() function.
#
# - Test row 2 (Sheet1 row 3) is removed; it was a duplicate test case.
# - The surviving test row (Sheet1 row 2) now reflects a re-run of the
#   find_crossing() backtest against ByBit data: new test number, new
#   symbol (BTCUSDT -> BTCUSD), a later "From" date, a corrected TP % and
#   the now-fractional SL % (the int/float bug), and the "Precision
#   Crossing" flag flipping to True.
# - Numeric columns F:J get a consistent 2-decimal-place number format
#   instead of the old mixed currency/plain styling, since they can now
#   hold true fractional values (e.g. 0.666...) instead of only ever
#   being whole-ish numbers.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("ListOfValues")

# --- Update the surviving data row (row 2) ---------------------------------
$ws1.Range("A2").Value = 11
$ws1.Range("B2").Value = "BTCUSD"
$ws1.Range("C2").Value = 44501   # 2021-11-01, stored as an Excel date serial
$ws1.Range("D2").Value = 44531   # 2021-12-01 (unchanged from before)
$ws1.Range("E2").Value = 30
$ws1.Range("F2").Value = 10000
$ws1.Range("G2").Value = 1
$ws1.Range("H2").Value = 2 / 3
$ws1.Range("I2").Value = -0.025
$ws1.Range("J2").Value = 0.075
$ws1.Range("K2").Value = "True"

# --- Remove the now-duplicate second test row (old row 3) ------------------
$ws1.Rows.Item(3).Delete()

# --- Re-style the numeric columns (Trade Amount .. Taker Fee %) ------------
# Previously a mix of a "$"-currency format (right aligned) and a plain
# general format; now a uniform 2-decimal numeric format, centered, with
# the header row keeping its bold font.
$dataRange = $ws1.Range("F2:J2")
$dataRange.NumberFormat = "0.00"
$dataRange.HorizontalAlignment = -4108
$dataRange.Font.Bold = $false

$headerRange = $ws1.Range("F1:J1")
$headerRange.NumberFormat = "0.00"
$headerRange.HorizontalAlignment = -4108
$headerRange.Font.Bold = $true

Write-Output "Sheet1 updated: row 3 deleted, row 2 refreshed, F:J re-formatted"
